# ----------------------------------------------------------------------
# Weekly CompStat refresh: new report week (18) covering 4/28-5/4/2025,
# with refreshed crime-complaint figures across precincts.
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/number and week-covering dates ---
$ws.Range("A8").Value = "Volume 32   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# --- Cells that flip from a numeric 0/value to the "no activity" text
#     marker (shared text "0"), matching the style of their text-typed
#     neighbours in the same row. ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Cells that flip the other way: from the text marker to an actual
#     numeric figure now that data is available. ---
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = 0
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("C31").Value = 1

# --- Refreshed weekly / 28-day / YTD / historical figures ---
$ws.Range("M15").Value = 25
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 53
$ws.Range("K16").Value = -45.283018867924
$ws.Range("L16").Value = -51.666666666666
$ws.Range("M16").Value = -35.555555555555
$ws.Range("N16").Value = -89.930555555555
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 9.090909090909
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 5.714285714285
$ws.Range("L17").Value = -36.206896551724
$ws.Range("M17").Value = 27.586206896551
$ws.Range("N17").Value = -51.315789473684
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = -36.585365853658
$ws.Range("L18").Value = -50.943396226415
$ws.Range("M18").Value = -24.637681159420
$ws.Range("N18").Value = -79.766536964980
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = -32.142857142857
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -12.987012987013
$ws.Range("I19").Value = 297
$ws.Range("J19").Value = 333
$ws.Range("K19").Value = -10.810810810810
$ws.Range("L19").Value = -28.605769230769
$ws.Range("M19").Value = -7.476635514018
$ws.Range("N19").Value = -60.766182298546
$ws.Range("D20").Value = 1
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = -68.75
$ws.Range("M20").Value = -58.333333333333
$ws.Range("N20").Value = -98.091603053435
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -13.513513513513
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -13.008130081300
$ws.Range("I21").Value = 425
$ws.Range("J21").Value = 520
$ws.Range("K21").Value = -18.269230769230
$ws.Range("L21").Value = -35.114503816793
$ws.Range("M21").Value = -11.458333333333
$ws.Range("N21").Value = -74.164133738601
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 19
$ws.Range("K22").Value = 18.75
$ws.Range("L22").Value = 5.555555555555
$ws.Range("M22").Value = -17.391304347826
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -40.476190476190
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 158
$ws.Range("H24").Value = -37.974683544303
$ws.Range("I24").Value = 489
$ws.Range("J24").Value = 586
$ws.Range("K24").Value = -16.552901023890
$ws.Range("L24").Value = -18.5
$ws.Range("M24").Value = -0.609756097560
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -64
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 119
$ws.Range("H25").Value = -55.462184873949
$ws.Range("I25").Value = 341
$ws.Range("J25").Value = 479
$ws.Range("K25").Value = -28.810020876826
$ws.Range("L25").Value = -22.5
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -3.030303030303
$ws.Range("I26").Value = 116
$ws.Range("J26").Value = 121
$ws.Range("K26").Value = -4.132231404958
$ws.Range("L26").Value = -17.142857142857
$ws.Range("M26").Value = 45
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = 12
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -60
$ws.Range("I31").Value = 6
$ws.Range("J31").Value = 8
$ws.Range("K31").Value = -25
$ws.Range("L31").Value = 500

$ws.Calculate()
